$d = $word.ActiveDocument

# --- Hunk 1: "jatekosok_szama" paragraph -------------------------------
# " 1 (alapertelmezes) vagy 2 lehet." -> " 1-4 lehet, 1 az alapertelmezes."
# (the italics on "alapertelmezes" are dropped as part of this edit)
$d.Content.Find.ClearFormatting()
$d.Content.Find.Replacement.ClearFormatting()
$d.Content.Find.Execute(
    ' 1 (alapértelmezés) vagy 2 lehet.',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ' 1-4 lehet, 1 az alapértelmezés.',
    2) | Out-Null

# --- Hunk 2: drop the "Ketjatekos uzemmodban ..." sentence --------------
$d.Content.Find.ClearFormatting()
$d.Content.Find.Replacement.ClearFormatting()
$d.Content.Find.Execute(
    ' A pálya széle szögletes és fehér, a kígyó kerek és színes karaktereket használ. Kétjátékos üzemmódban a két kígyó más színű.',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ' A pálya széle szögletes és fehér, a kígyó kerek és színes karaktereket használ.',
    2) | Out-Null

# --- Hunk 3: heading "Ketjatekos uzemmod" -> "Tobbjatekos uzemmod" ------
$d.Content.Find.ClearFormatting()
$d.Content.Find.Replacement.ClearFormatting()
$d.Content.Find.Execute(
    'Kétjátékos üzemmód',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'Többjátékos üzemmód',
    2) | Out-Null

# --- Hunk 4: rewrite the body paragraph of that section ------------------
$d.Content.Find.ClearFormatting()
$d.Content.Find.Replacement.ClearFormatting()
$d.Content.Find.Execute(
    'Ebben a módban két játékos versenyzik egymás ellen. A piros kígyót a WASD gombokkal, a másikat a nyilakkal lehet irányítani. Ha az egyik kígyó nekimegy vagy a falnak, vagy valamelyik kígyó testének, akkor ő vesztett, a győztes pontszáma pedig a két kígyó által összesen megevett gyümölcsök száma. A ranglista ehhez az üzemmódhoz különbözik az egyszemélyes játék ranglistájától. Matekfeladatot itt is meg lehet oldani, viszont a sikerességen felül a gyorsaság is számít. Aki eltalálta a jó megoldást (ezt nem a számítógép, hanem a két játékos kezeli egymás között), az kiválaszthatja, hogy melyik kígyó kezdjen a megszerzett pontok felével – eldöntheti, hogy őt segíti vagy hátráltatja-e a hosszúság. Hogyha sikertelen egy megoldás, akkor addig lehet újra próbálkozni, amíg valaki el nem találja, de fel is lehet adni. Ilyenkor mindkét kígyó alapmérettel kezd.',
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    'Ebben a módban több játékos versenyzik, megpróbálnak közösen minél több pontot gyűjteni. A zöld kígyót a nyilakkal, a kéket a WASD, a pirosat az IJKL, a sárgát pedig a TFGH gombokkal lehet irányítani. Egy gyümölcs megevésekor annyi pontot kap a csapat, ahány kígyó életben van még. Így érdemes minél több kígyót életben tartva enni az almákat. Ha az egyik kígyó nekimegy vagy a falnak, vagy valamelyik kígyó testének, akkor ő kiesik a játékból, a többiek nélküle játszanak tovább. A játékosok számától függően külön ranglistákban tároljuk az eredményeket. Matekfeladatot itt is meg lehet oldani, sikeres megoldás esetén a pontok felét szétosztja a kígyók között a program. Tehát ha 4 játékos 19 pontot gyűjtött, akkor 19/2=9 pontot kaphatnának. Mindenki 9/4=2 mezővel hosszabb kígyóval kezd, mint az alap méret, így végül is 8 ponttal kezdődik a játék. Hogyha sikertelen egy megoldás, akkor minden kígyó alapmérettel kezd. ',
    2) | Out-Null

# Re-apply italics to the four control-key hints in the rewritten paragraph.
function Set-ItalicOnce($searchText) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
    if ($found) {
        $rng.Italic = $true
    }
}

Set-ItalicOnce 'WASD'
Set-ItalicOnce 'IJKL'
Set-ItalicOnce 'TFGH '

Write-Output 'done'
